# Updates cryptos list prices/volume% (and swaps the Dai/Stacks rows 37-38)
# to reflect the latest scrape, mirroring the GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number ("105.50", "0.999", ...)
# are given an explicit Text number format first so Excel's COM layer
# stores them as strings (keeping trailing zeros, etc.) instead of
# silently coercing them into numeric values.

$ws.Range("D2").Value = "69.515.41"
$ws.Range("E2").Value = "  -0.83%  "

$ws.Range("D3").Value = "3.544.50"
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "197.33"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "584.53"
$ws.Range("E6").Value = "  -3.06%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  -2.28%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -1.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.632"
$ws.Range("E10").Value = "  -2.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "51.97"
$ws.Range("E11").Value = "  -3.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000288"
$ws.Range("E12").Value = "  -5.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.27"
$ws.Range("E13").Value = "  -2.98%  "

$ws.Range("D14").Value = "4.100.82"
$ws.Range("E14").Value = "  -2.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "666.39"
$ws.Range("E15").Value = "  +12.55%  "

$ws.Range("D16").Value = "69.592.29"
$ws.Range("E16").Value = "  -0.94%  "

$ws.Range("D17").Value = "3.545.48"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.42"
$ws.Range("E18").Value = "  -5.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.57"
$ws.Range("E19").Value = "  -3.32%  "

$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.970"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.46"
$ws.Range("E22").Value = "  +4.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "105.50"
$ws.Range("E24").Value = "  +3.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.39"
$ws.Range("E25").Value = "  -4.69%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.94"
$ws.Range("E26").Value = "  -3.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.22"
$ws.Range("E27").Value = "  -5.47%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("E28").Value = "  +1.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.53"
$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.46"
$ws.Range("E30").Value = "  -6.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.86"
$ws.Range("E31").Value = "  -3.82%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.89"
$ws.Range("E32").Value = "  -3.44%  "

$ws.Range("E33").Value = "  -4.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.96"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("D35").Value = "3.788.19"
$ws.Range("E35").Value = "  -3.54%  "

$ws.Range("D36").Value = "0.0₃0814"
$ws.Range("E36").Value = "  -8.59%  "

$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").Value = "  +4.59%  "

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "505.29"
$ws.Range("E39").Value = "  -4.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.93"
$ws.Range("E40").Value = "  -6.38%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  -4.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  +0.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.77"
$ws.Range("E43").Value = "  -6.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0455"
$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.91"
$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.137"
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("E49").Value = "  -3.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.79"
$ws.Range("E50").Value = "  +20.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  +62.62%  "
